# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $range = $ws.Range($cell)
    $range.Style = "Normal"
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell 'D2' '34.770.67'
Set-TextCell 'E2' '  -1.62%  '
Set-TextCell 'D3' '1.872.00'
Set-TextCell 'E3' '  -2.11%  '
Set-TextCell 'E4' '  -1.04%  '
Set-TextCell 'D5' '246.66'
Set-TextCell 'E5' '  -2.30%  '
Set-TextCell 'D6' '0.688'
Set-TextCell 'E6' '  -5.00%  '
Set-TextCell 'E7' '  -0.99%  '
Set-TextCell 'D8' '42.00'
Set-TextCell 'E8' '  +3.61%  '
Set-TextCell 'D9' '0.346'
Set-TextCell 'E9' '  -3.20%  '
Set-TextCell 'D10' '50.94'
Set-TextCell 'E10' '  -3.46%  '
Set-TextCell 'D11' '0.0735'
Set-TextCell 'E11' '  +0.35%  '
Set-TextCell 'D12' '0.0970'
Set-TextCell 'E12' '  -2.91%  '
Set-TextCell 'D13' '2.142.88'
Set-TextCell 'E13' '  -2.18%  '
Set-TextCell 'D14' '12.75'
Set-TextCell 'E14' '  +1.54%  '
Set-TextCell 'D15' '0.712'
Set-TextCell 'E15' '  -0.24%  '
Set-TextCell 'D16' '4.88'
Set-TextCell 'E16' '  -0.05%  '
Set-TextCell 'D17' '1.883.73'
Set-TextCell 'E17' '  -1.65%  '
Set-TextCell 'D18' '34.737.20'
Set-TextCell 'E18' '  -1.72%  '
Set-TextCell 'D19' '72.64'
Set-TextCell 'E19' '  -0.53%  '
Set-TextCell 'D20' '0.0₃0818'
Set-TextCell 'E20' '  -1.34%  '
Set-TextCell 'D21' '243.42'
Set-TextCell 'E21' '  +0.75%  '
Set-TextCell 'D22' '12.66'
Set-TextCell 'E22' '  -2.97%  '
Set-TextCell 'D23' '4.90'
Set-TextCell 'E23' '  -3.21%  '
Set-TextCell 'E24' '  -1.02%  '
Set-TextCell 'D25' '2.43'
Set-TextCell 'E25' '  +4.43%  '
Set-TextCell 'D26' '2.24'
Set-TextCell 'E26' '  -3.98%  '
Set-TextCell 'D27' '164.94'
Set-TextCell 'E27' '  -1.73%  '
Set-TextCell 'D28' '8.36'
Set-TextCell 'E28' '  -3.66%  '
Set-TextCell 'D29' '18.17'
Set-TextCell 'E29' '  -2.78%  '
Set-TextCell 'E30' '  -5.09%  '
Set-TextCell 'D31' '4.128.37'
Set-TextCell 'E31' '  -0.09%  '
Set-TextCell 'D32' '1.68'
Set-TextCell 'E32' '  +2.13%  '
Set-TextCell 'D33' '4.27'
Set-TextCell 'E33' '  -2.25%  '
Set-TextCell 'D34' '0.0575'
Set-TextCell 'E34' '  -0.69%  '
Set-TextCell 'D35' '4.14'
Set-TextCell 'E35' '  -2.35%  '
Set-TextCell 'D37' '0.826'
Set-TextCell 'E37' '  -9.31%  '
Set-TextCell 'D38' '1.98'
Set-TextCell 'E38' '  -2.25%  '
Set-TextCell 'D39' '1.53'
Set-TextCell 'E39' '  -22.89%  '
Set-TextCell 'D40' '97.58'
Set-TextCell 'E40' '  -1.61%  '
Set-TextCell 'D41' '16.84'
Set-TextCell 'E41' '  -3.65%  '
Set-TextCell 'D42' '0.0657'
Set-TextCell 'E42' '  +1.00%  '
Set-TextCell 'E43' '  -0.01%  '
Set-TextCell 'D44' '1.08'
Set-TextCell 'E44' '  -4.92%  '
Set-TextCell 'D45' '1.281.16'
Set-TextCell 'E45' '  -4.97%  '
Set-TextCell 'D46' '2.32'
Set-TextCell 'E46' '  -6.46%  '
Set-TextCell 'B47' 'Cronos'
Set-TextCell 'C47' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D47' '0.0788'
Set-TextCell 'E47' '  +8.42%  '
Set-TextCell 'B48' 'HuobiToken'
Set-TextCell 'C48' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D48' '2.40'
Set-TextCell 'E48' '  -0.93%  '
Set-TextCell 'E49' '  -1.93%  '
Set-TextCell 'D50' '12.07'
Set-TextCell 'E50' '  +4.81%  '
Set-TextCell 'E51' '  -4.72%  '
